$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "2024-05-11 11:29:48"
$ws.Range("B5").Value = 0.0004
